# "Add files via upload" — refresh the title slide text and the
# auto-populated "last saved" date shown on the master/layout date
# placeholders.

$p = $ppt.ActivePresentation

# --- Slide 1: replace the draft title with the real deck title ---
$s1 = $p.Slides.Item(1)
$titleShape = $s1.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Auction DataSet"

# (subtitle "File created on: 4/9/22 6:15:43 PM EDT" is left untouched)

# --- Refresh the date placeholder text everywhere it appears ---
$newDate = "4/9/22"

$master = $p.SlideMaster

# Slide master's own Date Placeholder
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16) {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

# Every slide layout's Date Placeholder
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}
